# Applies the updated FFXIV Leviathan market/profit figures scraped by the
# scheduled runner onto the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Table_*" sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 936.6875
$ws.Range("I41").Value = 1403.5555
$ws.Range("J41").Value = 336.42856
$ws.Range("K41").Value = 1403.5555
$ws.Range("L41").Value = 336.42856
$ws.Range("M41").Value = -963.5554999999999
$ws.Range("N41").Value = -1216.42856
# Row 62
$ws.Range("H62").Value = 4355.091
$ws.Range("I62").Value = 3977.4707
$ws.Range("J62").Value = 5639
$ws.Range("K62").Value = 3977.4707
$ws.Range("L62").Value = 5639
$ws.Range("M62").Value = -3353.4707
$ws.Range("N62").Value = -6887
# Row 65
$ws.Range("H65").Value = 4355.091
$ws.Range("I65").Value = 3977.4707
$ws.Range("J65").Value = 5639
$ws.Range("K65").Value = 19887.3535
$ws.Range("L65").Value = 28195
$ws.Range("M65").Value = -16767.3535
$ws.Range("N65").Value = -34435
# Row 100
$ws.Range("H100").Value = 7390.25
$ws.Range("I100").Value = 8031.364
$ws.Range("K100").Value = 8031.364
$ws.Range("M100").Value = -7490.364
# Row 132
$ws.Range("H132").Value = 3046.7646
$ws.Range("I132").Value = 2800.037
$ws.Range("J132").Value = 3998.4285
$ws.Range("K132").Value = 8400.110999999999
$ws.Range("L132").Value = 11995.2855
$ws.Range("M132").Value = -5870.110999999999
$ws.Range("N132").Value = -17055.2855
# Row 133
$ws.Range("H133").Value = 96963.336
$ws.Range("J133").Value = 96963.336
$ws.Range("L133").Value = 96963.336
$ws.Range("N133").Value = -107083.336
# Row 135
$ws.Range("H135").Value = 126377.375
$ws.Range("I135").Value = 1437
$ws.Range("J135").Value = 251317.75
$ws.Range("K135").Value = 12933
$ws.Range("L135").Value = 2261859.75
$ws.Range("M135").Value = -10398
$ws.Range("N135").Value = -2266929.75
# Row 139
$ws.Range("H139").Value = 219999
$ws.Range("J139").Value = 219999
$ws.Range("L139").Value = 219999
$ws.Range("N139").Value = -230279
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 2262.1914
$ws.Range("I74").Value = 1989.079
$ws.Range("J74").Value = 3415.3333
$ws.Range("K74").Value = 1989.079
$ws.Range("L74").Value = 3415.3333
$ws.Range("M74").Value = -1115.079
$ws.Range("N74").Value = -5163.3333
# Row 77
$ws.Range("H77").Value = 2262.1914
$ws.Range("I77").Value = 1989.079
$ws.Range("J77").Value = 3415.3333
$ws.Range("K77").Value = 9945.395
$ws.Range("L77").Value = 17076.6665
$ws.Range("M77").Value = -5577.395
$ws.Range("N77").Value = -25812.6665
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
# Row 19
$ws.Range("H19").Value = 1000000
$ws.Range("I19").Value = 1000000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1000000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -999827
$ws.Range("N19").ClearContents()
# Row 105
$ws.Range("H105").Value = 1871
$ws.Range("I105").Value = 2013.25
$ws.Range("K105").Value = 2013.25
$ws.Range("M105").Value = -266.25
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 200459.8
$ws.Range("I2").Value = 324.75
$ws.Range("J2").Value = 1001000
$ws.Range("K2").Value = 324.75
$ws.Range("L2").Value = 1001000
$ws.Range("M2").Value = -211.75
$ws.Range("N2").Value = -1001226
# Row 16
$ws.Range("H16").Value = 2895.5
$ws.Range("I16").Value = 3065.6428
$ws.Range("J16").Value = 2300
$ws.Range("K16").Value = 3065.6428
$ws.Range("L16").Value = 2300
$ws.Range("M16").Value = -2778.6428
$ws.Range("N16").Value = -2874
# Row 113
$ws.Range("H113").Value = 2895.5
$ws.Range("I113").Value = 3065.6428
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 3065.6428
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = -895.6428000000001
$ws.Range("N113").Value = -6640
# Row 122
$ws.Range("H122").Value = 176998.83
$ws.Range("I122").Value = 204798.8
$ws.Range("J122").Value = 37999
$ws.Range("K122").Value = 614396.3999999999
$ws.Range("L122").Value = 113997
$ws.Range("M122").Value = -611946.3999999999
$ws.Range("N122").Value = -118897
# Row 141
$ws.Range("H141").Value = 292499
$ws.Range("J141").Value = 292499
$ws.Range("L141").Value = 292499
$ws.Range("N141").Value = -302859
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 20784222
$ws.Range("I4").Value = 3125275
$ws.Range("J4").Value = 56102116
$ws.Range("K4").Value = 9375825
$ws.Range("L4").Value = 168306348
$ws.Range("M4").Value = -9375713
$ws.Range("N4").Value = -168306572
# Row 5
$ws.Range("H5").Value = 841.35486
$ws.Range("I5").Value = 796.9167
$ws.Range("J5").Value = 993.7143
$ws.Range("K5").Value = 2390.7501
$ws.Range("L5").Value = 2981.1429
$ws.Range("M5").Value = -2278.7501
$ws.Range("N5").Value = -3205.1429
# Row 7
$ws.Range("H7").Value = 58
$ws.Range("I7").Value = 44
$ws.Range("K7").Value = 132
$ws.Range("M7").Value = -20
# Row 132
$ws.Range("H132").Value = 1951.8889
$ws.Range("J132").Value = 2363.3635
$ws.Range("L132").Value = 21270.2715
$ws.Range("N132").Value = -26330.2715
# Row 133
$ws.Range("H133").Value = 9230.23
$ws.Range("I133").Value = 4498.25
$ws.Range("J133").Value = 11333.333
$ws.Range("K133").Value = 13494.75
$ws.Range("L133").Value = 33999.999
$ws.Range("M133").Value = -8434.75
$ws.Range("N133").Value = -44119.999
# Row 134
$ws.Range("H134").Value = 2746.6667
$ws.Range("I134").Value = 2746.6667
$ws.Range("K134").Value = 8240.000100000001
$ws.Range("M134").Value = -3170.000100000001
# Row 135
$ws.Range("H135").Value = 841.35486
$ws.Range("I135").Value = 796.9167
$ws.Range("J135").Value = 993.7143
$ws.Range("K135").Value = 7172.2503
$ws.Range("L135").Value = 8943.4287
$ws.Range("M135").Value = -4637.2503
$ws.Range("N135").Value = -14013.4287
$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 1313.1428
$ws.Range("I31").Value = 976.5
$ws.Range("J31").Value = 3333
$ws.Range("K31").Value = 976.5
$ws.Range("L31").Value = 3333
$ws.Range("M31").Value = -684.5
$ws.Range("N31").Value = -3917
# Row 37
$ws.Range("H37").Value = 1313.1428
$ws.Range("I37").Value = 976.5
$ws.Range("J37").Value = 3333
$ws.Range("K37").Value = 976.5
$ws.Range("L37").Value = 3333
$ws.Range("M37").Value = -699.5
$ws.Range("N37").Value = -3887
# Row 102
$ws.Range("H102").Value = 2310.394
$ws.Range("I102").Value = 2100.64
$ws.Range("K102").Value = 2100.64
$ws.Range("M102").Value = -478.6399999999999
# Row 122
$ws.Range("H122").Value = 1921.2
$ws.Range("I122").Value = 1800.8334
$ws.Range("K122").Value = 5402.5002
$ws.Range("M122").Value = -2952.5002
# Row 126
$ws.Range("H126").Value = 6159.1055
$ws.Range("I126").Value = 4997.5386
$ws.Range("J126").Value = 8675.833000000001
$ws.Range("K126").Value = 14992.6158
$ws.Range("L126").Value = 26027.499
$ws.Range("M126").Value = -12522.6158
$ws.Range("N126").Value = -30967.499
# Row 132
$ws.Range("H132").Value = 3637.8286
$ws.Range("I132").Value = 2856.4807
$ws.Range("K132").Value = 8569.4421
$ws.Range("M132").Value = -6039.4421
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 4836
$ws.Range("I22").Value = 4450
$ws.Range("K22").Value = 4450
$ws.Range("M22").Value = -4155
# Row 27
$ws.Range("H27").Value = 4836
$ws.Range("I27").Value = 4450
$ws.Range("K27").Value = 4450
$ws.Range("M27").Value = -4343
# Row 61
$ws.Range("H61").Value = 197178.17
$ws.Range("I61").Value = 240976.25
$ws.Range("K61").Value = 240976.25
$ws.Range("M61").Value = -240774.25
# Row 93
$ws.Range("H93").Value = 23973.312
$ws.Range("I93").Value = 3549
$ws.Range("K93").Value = 3549
$ws.Range("M93").Value = -2301
# Row 113
$ws.Range("H113").Value = 197178.17
$ws.Range("I113").Value = 240976.25
$ws.Range("K113").Value = 240976.25
$ws.Range("M113").Value = -238806.25
# Row 122
$ws.Range("H122").Value = 5182.5
$ws.Range("I122").Value = 4248.25
$ws.Range("J122").Value = 5649.625
$ws.Range("K122").Value = 12744.75
$ws.Range("L122").Value = 16948.875
$ws.Range("M122").Value = -10294.75
$ws.Range("N122").Value = -21848.875
# Row 132
$ws.Range("H132").Value = 3887.6
$ws.Range("I132").Value = 2221.5
$ws.Range("J132").Value = 4998.3335
$ws.Range("K132").Value = 6664.5
$ws.Range("L132").Value = 14995.0005
$ws.Range("M132").Value = -4134.5
$ws.Range("N132").Value = -20055.0005
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1930.8182
$ws.Range("I122").Value = 1930.8182
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5792.4546
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3342.4546
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 4031.9783
$ws.Range("I132").Value = 4604.7427
$ws.Range("K132").Value = 13814.2281
$ws.Range("M132").Value = -11284.2281
# Row 136
$ws.Range("H136").Value = 891.46344
$ws.Range("I136").Value = 567.74286
$ws.Range("J136").Value = 2779.8333
$ws.Range("K136").Value = 1703.22858
$ws.Range("L136").Value = 8339.499899999999
$ws.Range("M136").Value = 846.77142
$ws.Range("N136").Value = -13439.4999
# Row 139
$ws.Range("H139").Value = 135000
$ws.Range("J139").Value = 135000
$ws.Range("L139").Value = 135000
$ws.Range("N139").Value = -145280

Write-Host "Updated $($wb.Worksheets.Count) sheets with refreshed Leviathan profit data"
